$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 values - new match entry for Al-Taawon vs Al Fayha
$ws.Range("A7").Value = "23/10/2025"
$ws.Range("B7").Value = "Al Fayha"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Al-Taawon"
$ws.Range("F7").Value = "W"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 0.66
$ws.Range("L7").Value = 0.58
$ws.Range("M7").Value = 8
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 6
$ws.Range("P7").Value = 3
